# Update cryptocurrency price/volume data per the commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.311.96"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "'1.705.98"
$ws.Range("E3").Value = "  -1.07%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'224.04"
$ws.Range("E5").Value = "  -1.04%  "
$ws.Range("D6").Value = "'0.5307"
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D7").Value = "'1.004"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.2661"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "'0.06583"
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "'20.75"
$ws.Range("E10").Value = "  -4.56%  "
$ws.Range("D11").Value = "'0.07624"
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").Value = "'1.710.50"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "'1.940.60"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "'0.5775"
$ws.Range("E15").Value = "  -2.08%  "
$ws.Range("D16").Value = "'0.0₅8140"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").Value = "'67.61"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").Value = "'27.314.99"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "'215.45"
$ws.Range("E19").Value = "  -4.36%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "'4.615"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "'10.37"
$ws.Range("E22").Value = "  -3.48%  "
$ws.Range("D23").Value = "'5.958"
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'144.44"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("D26").Value = "'1.704"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "'0.1202"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "'7.214"
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("D29").Value = "'16.11"
$ws.Range("E29").Value = "  -3.91%  "
$ws.Range("D30").Value = "'0.05374"
$ws.Range("E30").Value = "  -3.80%  "
$ws.Range("D31").Value = "'1.288"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "'3.467"
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("D34").Value = "'1.646"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "'2.865"
$ws.Range("E35").Value = "  +1.77%  "
$ws.Range("D36").Value = "'2.418"
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'0.9459"
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("D38").Value = "'0.5797"
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("D39").Value = "'0.01626"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").Value = "'5.769"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("D41").Value = "'1.005"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "'0.8414"
$ws.Range("E42").Value = "  -2.25%  "
$ws.Range("D43").Value = "'1.040.50"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").Value = "'100.97"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("D45").Value = "'1.848.63"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").Value = "'0.0₈116"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("D47").Value = "'57.74"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").Value = "'0.4516"
$ws.Range("E48").Value = "  +1.91%  "
$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("E50").Value = "  -2.36%  "
$ws.Range("E51").Value = "  -1.08%  "
